# Update financial data values in the "company_list" sheet (error fixes to IFRS list)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

$ws.Range("D2").Value = 18750
$ws.Range("E2").Value = 183
$ws.Range("F2").Value = 183
$ws.Range("G2").Value = -402
$ws.Range("H2").Value = -575
$ws.Range("I2").Value = -695
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 28790
$ws.Range("L2").Value = 16458
$ws.Range("M2").Value = 12332
$ws.Range("N2").Value = 8311
$ws.Range("O2").Value = 4022
$ws.Range("P2").Value = 395
$ws.Range("Q2").Value = 40
$ws.Range("R2").Value = 138
$ws.Range("S2").Value = 275
$ws.Range("T2").Value = 623
$ws.Range("U2").Value = -583
$ws.Range("V2").Value = 6545
$ws.Range("W2").Value = 0.97
$ws.Range("X2").Value = -3.07
$ws.Range("Y2").Value = -8.16
$ws.Range("Z2").Value = -1.96
$ws.Range("AA2").Value = 133.45
$ws.Range("AB2").Value = 1983.2
$ws.Range("AC2").Value = -881
$ws.Range("AD2").Value = -5.57
$ws.Range("AE2").Value = 11732
$ws.Range("AF2").Value = 0.42
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 76400000
$ws.Range("D3").Value = 18835
$ws.Range("E3").Value = 568
$ws.Range("F3").Value = 568
$ws.Range("G3").Value = 395
$ws.Range("H3").Value = 166
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = 150
$ws.Range("K3").Value = 30113
$ws.Range("L3").Value = 17472
$ws.Range("M3").Value = 12641
$ws.Range("N3").Value = 8481
$ws.Range("O3").Value = 4160
$ws.Range("P3").Value = 395
$ws.Range("Q3").Value = 588
$ws.Range("R3").Value = -827
$ws.Range("S3").Value = 647
$ws.Range("T3").Value = 589
$ws.Range("U3").Value = -1
$ws.Range("V3").Value = 7562
$ws.Range("W3").Value = 3.02
$ws.Range("X3").Value = 0.88
$ws.Range("Y3").Value = 0.18
$ws.Range("Z3").Value = 0.5600000000000001
$ws.Range("AA3").Value = 138.21
$ws.Range("AB3").Value = 1988.16
$ws.Range("AC3").Value = 19
$ws.Range("AD3").Value = 266.46
$ws.Range("AE3").Value = 11973
$ws.Range("AF3").Value = 0.43
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 76400000
$ws.Range("D4").Value = 20593
$ws.Range("E4").Value = 971
$ws.Range("F4").Value = 971
$ws.Range("G4").Value = 390
$ws.Range("H4").Value = 18
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 17
$ws.Range("K4").Value = 34893
$ws.Range("L4").Value = 22280
$ws.Range("M4").Value = 12614
$ws.Range("N4").Value = 8272
$ws.Range("O4").Value = 4342
$ws.Range("P4").Value = 395
$ws.Range("Q4").Value = -1055
$ws.Range("R4").Value = -711
$ws.Range("S4").Value = 1917
$ws.Range("T4").Value = 571
$ws.Range("U4").Value = -1625
$ws.Range("V4").Value = 10273
$ws.Range("W4").Value = 4.71
$ws.Range("X4").Value = 0.09
$ws.Range("Y4").Value = 0.02
$ws.Range("Z4").Value = 0.06
$ws.Range("AA4").Value = 176.63
$ws.Range("AB4").Value = 1970.38
$ws.Range("AC4").Value = 2
$ws.Range("AD4").Value = 2826.46
$ws.Range("AE4").Value = 11678
$ws.Range("AF4").Value = 0.46
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 76400000
$ws.Range("D5").Value = 32664
$ws.Range("E5").Value = 3111
$ws.Range("F5").Value = 3111
$ws.Range("G5").Value = 2179
$ws.Range("H5").Value = 1233
$ws.Range("I5").Value = 903
$ws.Range("J5").Value = 330
$ws.Range("K5").Value = 47520
$ws.Range("L5").Value = 33166
$ws.Range("M5").Value = 14354
$ws.Range("N5").Value = 9569
$ws.Range("O5").Value = 4785
$ws.Range("P5").Value = 395
$ws.Range("Q5").Value = 270
$ws.Range("R5").Value = 1241
$ws.Range("S5").Value = -1133
$ws.Range("T5").Value = 1259
$ws.Range("U5").Value = -989
$ws.Range("V5").Value = 15438
$ws.Range("W5").Value = 9.529999999999999
$ws.Range("X5").Value = 3.77
$ws.Range("Y5").Value = 10.12
$ws.Range("Z5").Value = 2.99
$ws.Range("AA5").Value = 231.06
$ws.Range("AB5").Value = 2214.41
$ws.Range("AC5").Value = 1143
$ws.Range("AD5").Value = 9.01
$ws.Range("AE5").Value = 13509
$ws.Range("AF5").Value = 0.76
$ws.Range("AG5").Value = 90
$ws.Range("AH5").Value = 0.87
$ws.Range("AI5").Value = 7.08
$ws.Range("AJ5").Value = 76400000
$ws.Range("D6").Value = 38505
$ws.Range("E6").Value = 4635
$ws.Range("F6").Value = 4635
$ws.Range("G6").Value = 3771
$ws.Range("H6").Value = 2441
$ws.Range("I6").Value = 1894
$ws.Range("K6").Value = 49586
$ws.Range("L6").Value = 34763
$ws.Range("M6").Value = 14823
$ws.Range("N6").Value = 10144
$ws.Range("P6").Value = 395
$ws.Range("Q6").Value = 1150
$ws.Range("R6").Value = -1386
$ws.Range("S6").Value = -213
$ws.Range("T6").Value = 401
$ws.Range("U6").Value = 749
$ws.Range("V6").Value = 17287
$ws.Range("W6").Value = 12.04
$ws.Range("X6").Value = 6.34
$ws.Range("Y6").Value = 19.21
$ws.Range("Z6").Value = 5.03
$ws.Range("AA6").Value = 234.52
$ws.Range("AB6").Value = 2510.22
$ws.Range("AC6").Value = 2398
$ws.Range("AD6").Value = 4.73
$ws.Range("AE6").Value = 14320
$ws.Range("AF6").Value = 0.79
$ws.Range("AG6").Value = 125
$ws.Range("AH6").Value = 1.1
$ws.Range("AI6").Value = 4.68
$ws.Range("AJ6").Value = 76400000
$ws.Range("D7").Value = 38107
$ws.Range("E7").Value = 4103
$ws.Range("G7").Value = 3508
$ws.Range("H7").Value = 2149
$ws.Range("I7").Value = 1914
$ws.Range("K7").Value = 54514
$ws.Range("L7").Value = 37422
$ws.Range("M7").Value = 17091
$ws.Range("N7").Value = 12138
$ws.Range("P7").Value = 393
$ws.Range("Q7").Value = 2280
$ws.Range("R7").Value = -2306
$ws.Range("S7").Value = 1696
$ws.Range("T7").Value = 1259
$ws.Range("U7").Value = 750
$ws.Range("W7").Value = 10.77
$ws.Range("X7").Value = 5.64
$ws.Range("Y7").Value = 17.18
$ws.Range("Z7").Value = 4.13
$ws.Range("AA7").Value = 218.96
$ws.Range("AC7").Value = 2424
$ws.Range("AD7").Value = 5.69
$ws.Range("AE7").Value = 17136
$ws.Range("AF7").Value = 0.8100000000000001
$ws.Range("AG7").Value = 127
$ws.Range("AH7").Value = 0.92
$ws.Range("AI7").Value = 5.05
$ws.Range("D8").Value = 35375
$ws.Range("E8").Value = 4058
$ws.Range("G8").Value = 3700
$ws.Range("H8").Value = 2787
$ws.Range("I8").Value = 2428
$ws.Range("K8").Value = 54319
$ws.Range("L8").Value = 34663
$ws.Range("M8").Value = 19655
$ws.Range("N8").Value = 14306
$ws.Range("P8").Value = 393
$ws.Range("Q8").Value = 2882
$ws.Range("R8").Value = -658
$ws.Range("S8").Value = -2475
$ws.Range("T8").Value = 430
$ws.Range("U8").Value = 1755
$ws.Range("W8").Value = 11.47
$ws.Range("X8").Value = 7.88
$ws.Range("Y8").Value = 18.36
$ws.Range("Z8").Value = 5.12
$ws.Range("AA8").Value = 176.36
$ws.Range("AC8").Value = 3075
$ws.Range("AD8").Value = 4.49
$ws.Range("AE8").Value = 20195
$ws.Range("AF8").Value = 0.68
$ws.Range("AG8").Value = 127
$ws.Range("AH8").Value = 0.92
$ws.Range("AI8").Value = 3.99
$ws.Range("D9").Value = 38709
$ws.Range("E9").Value = 4380
$ws.Range("G9").Value = 4039
$ws.Range("H9").Value = 3031
$ws.Range("I9").Value = 2797
$ws.Range("K9").Value = 57832
$ws.Range("L9").Value = 35369
$ws.Range("M9").Value = 22463
$ws.Range("N9").Value = 16636
$ws.Range("P9").Value = 393
$ws.Range("Q9").Value = 2576
$ws.Range("R9").Value = -1128
$ws.Range("S9").Value = -983
$ws.Range("T9").Value = 704
$ws.Range("U9").Value = 1400
$ws.Range("W9").Value = 11.31
$ws.Range("X9").Value = 7.83
$ws.Range("Y9").Value = 18.08
$ws.Range("Z9").Value = 5.41
$ws.Range("AA9").Value = 157.45
$ws.Range("AC9").Value = 3542
$ws.Range("AD9").Value = 3.9
$ws.Range("AE9").Value = 23484
$ws.Range("AF9").Value = 0.59
$ws.Range("AG9").Value = 127
$ws.Range("AH9").Value = 0.92
$ws.Range("AI9").Value = 3.46
